$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$helper = $ws.Range("Z1")

# Step 1: fill + alignment from L9
$srcFill = $ws.Range("L9")
$srcFill.Copy()
$helper.PasteSpecial(-4122)

# Step 2: border left/right thick only (weight before linestyle to avoid transient style)
$left = $helper.Borders.Item(7)
$left.Weight = 4
$left.LineStyle = 1
$left.ColorIndex = 1

$right = $helper.Borders.Item(10)
$right.Weight = 4
$right.LineStyle = 1
$right.ColorIndex = 1

Write-Output "done"
